$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "User Stories Planejadas" (sheet1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("User Stories Planejadas")

# Remove the "Prioridade" value (column C) for the first two rows.
$ws1.Range("C2").ClearContents()
$ws1.Range("C3").ClearContents()

# Update User Story Point values (column D) and add Sprint numbers (column E).
$ws1.Range("D4").Value = 13
$ws1.Range("E4").Value = 2

$ws1.Range("D5").Value = 40
$ws1.Range("E5").Value = 2

$ws1.Range("D6").Value = 40
$ws1.Range("E6").Value = 2

$ws1.Range("D7").Value = 20
$ws1.Range("E7").Value = 3

$ws1.Range("D8").Value = 100
$ws1.Range("E8").Value = 3

$ws1.Range("D9").Value = 40
$ws1.Range("E9").Value = 4

# Row 10 ("A10" = 9) is no longer part of the table - clear it out.
$ws1.Range("A10").ClearContents()

$ws1.Range("D10").Select()

# ---------------------------------------------------------------------------
# Sheet "User Stories Realizadas" (sheet2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("User Stories Realizadas")

$ws2.Range("D3").Value = 40
$ws2.Range("D4").Value = 40
$ws2.Range("D5").Value = 40
$ws2.Range("D6").Value = 3
$ws2.Range("D7").Value = 40
$ws2.Range("D8").Value = 100
$ws2.Range("D9").Value = 20

$ws2.Range("D8").Select()

# ---------------------------------------------------------------------------
# Sheet "Grafico - Use Storie" (sheet3)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Grafico - Use Storie")
$ws3.Range("P14").Select()

$wb.Save()
